# Update column G ("K") values per regenerated save_data (stat computed as
# actual strikeouts instead of legacy "Strike#" figure).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 6
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 1
    15 = 0
    16 = 3
    17 = 1
    18 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
